$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the generated-on timestamp in the report header
$ws.Range("A1").Value = "Reporte generado el 18/05/2025 a las 20:52"

# Update inventory rows 10-17 with corrected data (special characters bug fix).
# Row 18's data ("Vas A Caer Chuponsito") moves up into row 17, and the
# now-duplicate trailing row 18 is removed from the sheet.
$ws.Range("B10").Value = "Juego Sala - 123"
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "unidad"
$ws.Range("F10").Value = "Comedor"

$ws.Range("B11").Value = "Juego Sala- 123"
$ws.Range("C11").Value = 12
$ws.Range("D11").Value = 25.55
$ws.Range("E11").Value = "juego"
$ws.Range("F11").Value = "Sala"

$ws.Range("B12").Value = "...Afsbhehuhfhfdjkdsfjksjksdfjkfdjsfdhjsfd"
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 111
$ws.Range("E12").Value = "kit"
$ws.Range("F12").Value = "Comedor"

$ws.Range("B13").Value = "pan"
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = "unidad"
$ws.Range("F13").Value = "Comedor"

$ws.Range("B14").Value = "Hola"
$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 123
$ws.Range("E14").Value = "unidad"
$ws.Range("F14").Value = "Comedor"

$ws.Range("B15").Value = "Hola 123"
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = "unidad"
$ws.Range("F15").Value = "Sala"

$ws.Range("B16").Value = "La Potona-"
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = "juego"
$ws.Range("F16").Value = "Oficina"

$ws.Range("B17").Value = "Vas A Cáer Chuponsito"
$ws.Range("C17").Value = 1212
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = "par"
$ws.Range("F17").Value = "Sala"

# Remove the now-duplicated last row (18)
$ws.Rows.Item(18).Delete()
